# Prepare for importing formula expressions to the model.
#
# On the "Formula" worksheet:
#   - Fill column B (rows 5-13) with the values 10..18.
#   - Fill column A (rows 5-13) with =CONCATENATE("MW","-",B<row>).
#     Row 5 is entered on its own; rows 6-13 are entered as one multi-cell
#     Range.Formula assignment so Excel records them as a single shared
#     formula group (t="shared" ref="A6:A13" si="...") exactly like a
#     fill-down of the formula from A6 through A13 would.
#   - Move the active selection from G14 to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formula")

for ($row = 5; $row -le 13; $row++) {
    $ws.Range("B$row").Value = $row + 5
}

$ws.Range("A5").Formula = '=CONCATENATE("MW","-",B5)'
$ws.Range("A6:A13").Formula = '=CONCATENATE("MW","-",B6)'

$ws.Activate()
$ws.Range("B7").Select()
